$wb = $excel.ActiveWorkbook

# --- Sheet1 = ASKARI -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ASKARI")

# Row 71: new E71 + changed G/I/K + new M71
$ws1.Range("E71").Value = 70000000
$ws1.Range("G71").Value = 300000000
$ws1.Range("I71").Value = 91506.84931506851
$ws1.Range("K71").Value = 1444454.794520548
$ws1.Range("M71").Value = 10.74

# Rows 72-84: changed G/I/K
$ws1.Range("G72").Value = 300000000
$ws1.Range("I72").Value = 91506.84931506851
$ws1.Range("K72").Value = 1535961.643835616

$ws1.Range("G73").Value = 300000000
$ws1.Range("I73").Value = 91506.84931506851
$ws1.Range("K73").Value = 1627468.493150685

$ws1.Range("G74").Value = 300000000
$ws1.Range("I74").Value = 91506.84931506851
$ws1.Range("K74").Value = 1718975.342465753

$ws1.Range("G75").Value = 300000000
$ws1.Range("I75").Value = 91506.84931506851
$ws1.Range("K75").Value = 1810482.191780821

$ws1.Range("G76").Value = 300000000
$ws1.Range("I76").Value = 91506.84931506851
$ws1.Range("K76").Value = 1901989.04109589

$ws1.Range("G77").Value = 300000000
$ws1.Range("I77").Value = 91506.84931506851
$ws1.Range("K77").Value = 1993495.890410958

$ws1.Range("G78").Value = 300000000
$ws1.Range("I78").Value = 91506.84931506851
$ws1.Range("K78").Value = 2085002.739726027

$ws1.Range("G79").Value = 300000000
$ws1.Range("I79").Value = 91506.84931506851
$ws1.Range("K79").Value = 2176509.589041095

$ws1.Range("G80").Value = 300000000
$ws1.Range("I80").Value = 91506.84931506851
$ws1.Range("K80").Value = 2268016.438356164

$ws1.Range("G81").Value = 300000000
$ws1.Range("I81").Value = 91506.84931506851
$ws1.Range("K81").Value = 2359523.287671233

$ws1.Range("G82").Value = 300000000
$ws1.Range("I82").Value = 91506.84931506851
$ws1.Range("K82").Value = 2451030.136986301

$ws1.Range("G83").Value = 300000000
$ws1.Range("I83").Value = 91506.84931506851
$ws1.Range("K83").Value = 2542536.98630137

$ws1.Range("G84").Value = 300000000
$ws1.Range("I84").Value = 91506.84931506851
$ws1.Range("K84").Value = 2634043.835616439

# New rows 85-88
$ws1.Range("B85").Value = 83
$ws1.Range("C85").Value = "14/02/2022"
$ws1.Range("G85").Value = 300000000
$ws1.Range("I85").Value = 91506.84931506851
$ws1.Range("K85").Value = 2725550.684931507

$ws1.Range("B86").Value = 84
$ws1.Range("C86").Value = "15/02/2022"
$ws1.Range("G86").Value = 300000000
$ws1.Range("I86").Value = 91506.84931506851
$ws1.Range("K86").Value = 2817057.534246576

$ws1.Range("B87").Value = 85
$ws1.Range("C87").Value = "16/02/2022"
$ws1.Range("G87").Value = 300000000
$ws1.Range("I87").Value = 91506.84931506851
$ws1.Range("K87").Value = 2908564.383561645

$ws1.Range("B88").Value = 86
$ws1.Range("C88").Value = "17/02/2022"
$ws1.Range("G88").Value = 300000000
$ws1.Range("I88").Value = 91506.84931506851
$ws1.Range("K88").Value = 3000071.232876713

# --- Sheet2 = DIBL ----------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DIBL")

$ws2.Range("B51").Value = 49
$ws2.Range("C51").Value = "14/02/2022"
$ws2.Range("G51").Value = 45000000
$ws2.Range("I51").Value = 16397.2602739726
$ws2.Range("K51").Value = 803465.753424658

$ws2.Range("B52").Value = 50
$ws2.Range("C52").Value = "15/02/2022"
$ws2.Range("G52").Value = 45000000
$ws2.Range("I52").Value = 16397.2602739726
$ws2.Range("K52").Value = 819863.0136986306

$ws2.Range("B53").Value = 51
$ws2.Range("C53").Value = "16/02/2022"
$ws2.Range("G53").Value = 45000000
$ws2.Range("I53").Value = 16397.2602739726
$ws2.Range("K53").Value = 836260.2739726033

$ws2.Range("B54").Value = 52
$ws2.Range("C54").Value = "17/02/2022"
$ws2.Range("G54").Value = 45000000
$ws2.Range("I54").Value = 16397.2602739726
$ws2.Range("K54").Value = 852657.5342465759

# --- Sheet3 = HBL -------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("HBL")

$ws3.Range("B79").Value = 77
$ws3.Range("C79").Value = "14/02/2022"
$ws3.Range("G79").Value = 28600000
$ws3.Range("I79").Value = 9285.205479452054
$ws3.Range("K79").Value = 688844.7123287665

$ws3.Range("B80").Value = 78
$ws3.Range("C80").Value = "15/02/2022"
$ws3.Range("G80").Value = 28600000
$ws3.Range("I80").Value = 9285.205479452054
$ws3.Range("K80").Value = 698129.9178082185

$ws3.Range("B81").Value = 79
$ws3.Range("C81").Value = "16/02/2022"
$ws3.Range("G81").Value = 28600000
$ws3.Range("I81").Value = 9285.205479452054
$ws3.Range("K81").Value = 707415.1232876705

$ws3.Range("B82").Value = 80
$ws3.Range("C82").Value = "17/02/2022"
$ws3.Range("G82").Value = 28600000
$ws3.Range("I82").Value = 9285.205479452054
$ws3.Range("K82").Value = 716700.3287671226

# --- Sheet4 = MEEZAN -----------------------------------------------------------
$ws4 = $wb.Worksheets.Item("MEEZAN")

$ws4.Range("B36").Value = 34
$ws4.Range("C36").Value = "14/02/2022"
$ws4.Range("G36").Value = 209549543.44
$ws4.Range("I36").Value = 42480.82812324383
$ws4.Range("K36").Value = 1133216.542315135

$ws4.Range("B37").Value = 35
$ws4.Range("C37").Value = "15/02/2022"
$ws4.Range("G37").Value = 209549543.44
$ws4.Range("I37").Value = 42480.82812324383
$ws4.Range("K37").Value = 1175697.370438378

$ws4.Range("B38").Value = 36
$ws4.Range("C38").Value = "16/02/2022"
$ws4.Range("G38").Value = 209549543.44
$ws4.Range("I38").Value = 42480.82812324383
$ws4.Range("K38").Value = 1218178.198561622

$ws4.Range("B39").Value = 37
$ws4.Range("C39").Value = "17/02/2022"
$ws4.Range("G39").Value = 209549543.44
$ws4.Range("I39").Value = 42480.82812324383
$ws4.Range("K39").Value = 1260659.026684866
